# Minor tweaks to log messages and comments
# - Adds a "CI in mass histories" / "diff" comparison to the OutputData sheet
#   (columns L and M), comparing the computed 95% CI (column K) against the
#   CI value recorded in each weight's mass history, plus the signed
#   difference between the two.
# - Leaves cursor/selection in a few cells as a side effect of the edit.

$wb = $excel.ActiveWorkbook

$wsOut = $wb.Worksheets.Item("OutputData")
$wsOut.Activate()

# Headers for the two new columns.
$wsOut.Range("L1").Value = "CI in mass histories"
$wsOut.Range("M1").Value = "diff"

# Values recorded in each weight's mass history file (column L), and the
# difference against the freshly-computed 95% CI in column K (column M).
$ciHistory = @{
    2  = 0.84191399999999994
    3  = 0.84255400000000003
    4  = 0.85037368837708471
    5  = 64.18602605707224
    6  = 0.84190599999999993
    7  = 28.814466884851551
    8  = 0.84190799999999999
    9  = 43.836285342311619
    10 = 20.26923914929991
    11 = 20.842508177843392
}

foreach ($row in 2..11) {
    $wsOut.Range("L$row").Value = $ciHistory[$row]
    $wsOut.Range("M$row").Formula = "=K$row-L$row"
}

# Rows 12 and 13 have no computed 95% CI (no history match), so only the
# "diff" column is touched, left blank, matching the style of the column
# above it.
$wsOut.Range("M12").Value = $null
$wsOut.Range("M13").Value = $null

$wsOut.Range("N8").Select()

# Cursor was left elsewhere on the other two sheets touched this session.
$wsInput = $wb.Worksheets.Item("InputData")
$wsInput.Activate()
$wsInput.Range("D21").Select()

$wsAdmin = $wb.Worksheets.Item("Admin")
$wsAdmin.Activate()
$wsAdmin.Range("H33").Select()

Write-Output "Applied CI-in-mass-histories diff columns to OutputData"
